# Insert 4 new data rows (a new reporting date, serial 44644 = 2022-03-24) at
# the top of the existing data block (row 584), pushing the previously-existing
# rows 584:675 down to 588:679. This mirrors the source diff, which shows the
# dimension growing from A1:T675 to A1:T679 and every row from 588 on being
# identical to the old row minus 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 584:675 down by 4 rows, leaving 584:587 empty (format
# copied down from row 583, which matches the original file's look).
$ws.Rows("584:587").Insert()

# Static column values shared by every record in this sheet/table.
$mercadoId  = 7
$mercado    = "Terminal Hortofrutícola Agro Chillán"
$region     = "Ñuble"
$codreg     = 16
$tipo       = "Fruta"
$productoId = 100104
$producto   = "Frutos de pepita"
$categoriaId = 100104002
$categoria  = "Manzana"
$unidad     = "$/caja 16 kilos empedrada"
$origen     = "Provincia de Curicó"
$kgUnidad   = 16
$fecha      = 44644

$newRows = @(
    @{ Variedad = "Fuji royal";    Calidad = "Especial"; Volumen = 60;  PMin = 12000; PMax = 12000; PProm = 12000; PKg = 750 },
    @{ Variedad = "Fuji royal";    Calidad = "Primera";  Volumen = 120; PMin = 10000; PMax = 11000; PProm = 10500; PKg = 656 },
    @{ Variedad = "Granny Smith";  Calidad = "Especial"; Volumen = 80;  PMin = 11000; PMax = 11000; PProm = 11000; PKg = 688 },
    @{ Variedad = "Granny Smith";  Calidad = "Primera";  Volumen = 160; PMin = 9000;  PMax = 10000; PProm = 9500;  PKg = 594 }
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 584 + $i
    $d = $newRows[$i]

    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $d.Variedad
    $ws.Cells.Item($r, 12).Value = $d.Calidad
    $ws.Cells.Item($r, 13).Value = $d.Volumen
    $ws.Cells.Item($r, 14).Value = $d.PMin
    $ws.Cells.Item($r, 15).Value = $d.PMax
    $ws.Cells.Item($r, 16).Value = $d.PProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $d.PKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
